# Append 39 new rows (r=1085..1123) of SOLUSDT 4h OHLCV candles to Sheet1,
# extending the used range from A1:F1084 to A1:F1123.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (including the date-time style on column A) from the last existing
# data row (1084) down across the new rows 1085:1123, so every appended row matches the
# style already used for all prior data rows (bold/bordered/centered date cell in col A).
$ws.Range("A1084:F1084").Copy()
$ws.Range("A1085:F1123").PasteSpecial(-4122)

# timestamp, open, high, low, close, volume
$data = @(
    @(45534.5, 138.65, 140.73, 132.49, 133.32, 1341211.748),
    @(45534.66666666666, 133.33, 137.91, 132.04, 135.98, 1131431.907),
    @(45534.83333333334, 135.98, 138.48, 135.65, 137.99, 320593.048),
    @(45535, 137.99, 139.23, 137.5, 138.29, 264961.297),
    @(45535.16666666666, 138.29, 138.77, 137.38, 137.59, 130327.815),
    @(45535.33333333334, 137.59, 137.85, 135.79, 137.01, 300188.45),
    @(45535.5, 137.02, 137.29, 135.57, 136.63, 177030.014),
    @(45535.66666666666, 136.63, 136.7, 134.49, 134.97, 178614.814),
    @(45535.83333333334, 134.97, 135.55, 134.51, 135.47, 86965.66099999999),
    @(45536, 135.35, 135.62, 132.88, 134.07, 279306.294),
    @(45536.16666666666, 134.07, 134.09, 131.09, 132.81, 397081.784),
    @(45536.33333333334, 132.81, 132.84, 130.18, 131.91, 436844.474),
    @(45536.5, 131.91, 134.87, 130, 134.29, 529518.341),
    @(45536.66666666666, 134.29, 135.92, 131.87, 135.23, 363573.128),
    @(45536.83333333334, 135.24, 136.15, 127.66, 128.5, 922120.72),
    @(45537, 128.5, 130.7, 128.27, 130.21, 567146.12),
    @(45537.16666666666, 130.21, 130.37, 126.8, 127.66, 639972.608),
    @(45537.33333333334, 127.65, 132.93, 127.22, 131.9, 785422.107),
    @(45537.5, 131.83, 133.45, 130, 132.54, 511501.253),
    @(45537.66666666666, 132.54, 133.87, 131.71, 133.05, 286378.41),
    @(45537.83333333334, 133.04, 133.5, 132.81, 133.45, 14352.089),
    @(45538, 134.95, 138.13, 134.21, 135.51, 401897.702),
    @(45538.16666666666, 135.5, 135.89, 133.54, 134.6, 337314.591),
    @(45538.33333333334, 134.6, 134.64, 132.77, 134.36, 279259.928),
    @(45538.5, 134.36, 135.5, 129.31, 129.8, 667203.267),
    @(45538.66666666666, 129.81, 132.21, 129.31, 130.95, 374791.042),
    @(45538.83333333334, 130.96, 131.61, 127.15, 127.54, 566915.395),
    @(45539, 127.54, 128.89, 122.6, 127.7, 1200406.418),
    @(45539.16666666666, 127.7, 131.44, 126.55, 130.3, 505116.797),
    @(45539.33333333334, 130.29, 130.71, 128.11, 128.43, 398907.541),
    @(45539.5, 128.43, 128.72, 128.4, 128.65, 3550.511),
    @(45539.66666666666, 133.82, 134.98, 132.58, 132.8, 459187.569),
    @(45539.83333333334, 133.19, 134.66, 132.89, 133.64, 203183.493),
    @(45540, 133.64, 134.38, 131.28, 131.6, 292498.457),
    @(45540.16666666666, 131.6, 134.26, 130.84, 133.1, 373792.043),
    @(45540.33333333334, 133.1, 133.87, 131.69, 132.83, 315818.006),
    @(45540.5, 132.83, 135.58, 130.16, 130.61, 905450.0330000001),
    @(45540.66666666666, 130.61, 131.14, 127.85, 128.44, 622556.003),
    @(45540.83333333334, 128.44, 129.91, 128.13, 128.92, 185115.496)
)

$r = 1085
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
